$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.161.08'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.506.10'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.89'
$ws.Range('E5').Value = '  -0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.86'
$ws.Range('E6').Value = '  +3.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.610'
$ws.Range('E7').Value = '  -0.93%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.500.99'
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.193'
$ws.Range('E10').Value = '  -0.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.18'
$ws.Range('E11').Value = '  +7.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.582'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '46.14'
$ws.Range('E13').Value = '  -1.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000275'
$ws.Range('E14').Value = '  -0.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.060.15'
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '609.80'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.503.97'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.259.37'
$ws.Range('E19').Value = '  +1.26%  '
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.31'
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.00'
$ws.Range('E23').Value = '  -11.50%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '98.44'
$ws.Range('E24').Value = '  +3.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '15.53'
$ws.Range('E25').Value = '  -1.41%  '
$ws.Range('E26').Value = '  -3.43%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('E28').Value = '  -1.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.83'
$ws.Range('E29').Value = '  +2.20%  '
$ws.Range('E31').Value = '  -4.47%  '
$ws.Range('E32').Value = '  -3.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.27'
$ws.Range('E33').Value = '  -3.56%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.81'
$ws.Range('E34').Value = '  -0.81%  '
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '627.39'
$ws.Range('E35').Value = '  +13.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0992'
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.53'
$ws.Range('E37').Value = '  +1.98%  '
$ws.Range('B38').Value = 'Cosmos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.73'
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0475'
$ws.Range('E39').Value = '  +6.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '56.69'
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.142'
$ws.Range('E42').Value = '  +2.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.363.79'
$ws.Range('E43').Value = '  +1.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₃0733'
$ws.Range('E44').Value = '  +5.31%  '
$ws.Range('E45').Value = '  -5.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '32.13'
$ws.Range('E46').Value = '  -2.48%  '
$ws.Range('E47').Value = '  +0.80%  '
$ws.Range('E48').Value = '  -1.52%  '
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.48'
$ws.Range('E50').Value = '  -2.17%  '
